$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new column F (Brukernavn) between Navn (E) and Programmeringserfaring (old F)
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).ColumnWidth = 19.2

# Move "Navn" old values (std1..std28) into new Brukernavn column, and set Navn to student1..student28
for ($r = 2; $r -le 29; $r++) {
    $n = $r - 1
    $oldNavn = $ws.Cells.Item($r, 5).Value()
    if ($r -eq 2) {
        $ws.Cells.Item($r, 6).Value = "std1"
    } else {
        $ws.Cells.Item($r, 6).Value = $oldNavn
    }
    $ws.Cells.Item($r, 5).Value = "student$n"
}

$ws.Range("F1").Value = "Brukernavn"

# Rebuild the table so column metadata (names) resync correctly
$tbl.Unlist()
$newTbl = $ws.ListObjects.Add(1, $ws.Range("A1:I29"), $null, 1)
$newTbl.Name = "Table1"

# Update sheet view per diff
$ws.Range("F2:F29").Select()
